# Auto-generated edits applying the diff against Sheets/Bahamut_Profits.xlsx
# (values for columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR tables)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 965.4545000000001
$ws.Range("I19").Value = 634.2857
$ws.Range("J19").Value = 1545
$ws.Range("K19").Value = 634.2857
$ws.Range("L19").Value = 1545
$ws.Range("M19").Value = -459.2857
$ws.Range("N19").Value = -1895
$ws.Range("H88").Value = 795010.7
$ws.Range("J88").Value = 1221831.6
$ws.Range("L88").Value = 1221831.6
$ws.Range("N88").Value = -1222643.6
$ws.Range("H91").Value = 795010.7
$ws.Range("J91").Value = 1221831.6
$ws.Range("L91").Value = 1221831.6
$ws.Range("N91").Value = -1224639.6
$ws.Range("H106").Value = 3537.5
$ws.Range("I106").Value = 3339.8572
$ws.Range("J106").Value = 3998.6667
$ws.Range("K106").Value = 3339.8572
$ws.Range("L106").Value = 3998.6667
$ws.Range("M106").Value = -2708.8572
$ws.Range("N106").Value = -5260.6667
$ws.Range("H132").Value = 154989.3
$ws.Range("I132").Value = 1120.9828
$ws.Range("K132").Value = 3362.9484
$ws.Range("M132").Value = -832.9484000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1547.1875
$ws.Range("I45").Value = 1150.4546
$ws.Range("J45").Value = 2420
$ws.Range("K45").Value = 1150.4546
$ws.Range("L45").Value = 2420
$ws.Range("M45").Value = -773.4546
$ws.Range("N45").Value = -3174
$ws.Range("H61").Value = 1329.7826
$ws.Range("I61").Value = 1354.2354
$ws.Range("J61").Value = 1260.5
$ws.Range("K61").Value = 1354.2354
$ws.Range("L61").Value = 1260.5
$ws.Range("M61").Value = -1142.2354
$ws.Range("N61").Value = -1684.5
$ws.Range("H88").Value = 2587
$ws.Range("I88").Value = 2252
$ws.Range("J88").Value = 2773.111
$ws.Range("K88").Value = 2252
$ws.Range("L88").Value = 2773.111
$ws.Range("M88").Value = -1846
$ws.Range("N88").Value = -3585.111
$ws.Range("H91").Value = 2587
$ws.Range("I91").Value = 2252
$ws.Range("J91").Value = 2773.111
$ws.Range("K91").Value = 2252
$ws.Range("L91").Value = 2773.111
$ws.Range("M91").Value = -848
$ws.Range("N91").Value = -5581.111
$ws.Range("H136").Value = 1329.7826
$ws.Range("I136").Value = 1354.2354
$ws.Range("J136").Value = 1260.5
$ws.Range("K136").Value = 4062.7062
$ws.Range("L136").Value = 3781.5
$ws.Range("M136").Value = -1512.7062
$ws.Range("N136").Value = -8881.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 40000
$ws.Range("J6").Value = 40000
$ws.Range("L6").Value = 40000
$ws.Range("N6").Value = -40226
$ws.Range("H86").Value = 2209.7144
$ws.Range("I86").Value = 1941.4117
$ws.Range("J86").Value = 3350
$ws.Range("K86").Value = 1941.4117
$ws.Range("L86").Value = 3350
$ws.Range("M86").Value = -818.4117000000001
$ws.Range("N86").Value = -5596
$ws.Range("H89").Value = 2209.7144
$ws.Range("I89").Value = 1941.4117
$ws.Range("J89").Value = 3350
$ws.Range("K89").Value = 9707.058500000001
$ws.Range("L89").Value = 16750
$ws.Range("M89").Value = -4091.058500000001
$ws.Range("N89").Value = -27982
$ws.Range("H105").Value = 4458.5186
$ws.Range("I105").Value = 4493.5405
$ws.Range("J105").Value = 4382.294
$ws.Range("K105").Value = 4493.5405
$ws.Range("L105").Value = 4382.294
$ws.Range("M105").Value = -2746.5405
$ws.Range("N105").Value = -7876.294
$ws.Range("H114").Value = 25000
$ws.Range("I114").Value = 25000
$ws.Range("K114").Value = 25000
$ws.Range("M114").Value = -20661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1895.2821
$ws.Range("I31").Value = 1889.6216
$ws.Range("K31").Value = 1889.6216
$ws.Range("M31").Value = -1594.6216
$ws.Range("H34").Value = 1895.2821
$ws.Range("I34").Value = 1889.6216
$ws.Range("K34").Value = 1889.6216
$ws.Range("M34").Value = -1687.6216
$ws.Range("H62").Value = 6789.1787
$ws.Range("I62").Value = 7002.778
$ws.Range("J62").Value = 6404.7
$ws.Range("K62").Value = 7002.778
$ws.Range("L62").Value = 6404.7
$ws.Range("M62").Value = -6378.778
$ws.Range("N62").Value = -7652.7
$ws.Range("H65").Value = 6789.1787
$ws.Range("I65").Value = 7002.778
$ws.Range("J65").Value = 6404.7
$ws.Range("K65").Value = 35013.89
$ws.Range("L65").Value = 32023.5
$ws.Range("M65").Value = -31893.89
$ws.Range("N65").Value = -38263.5
$ws.Range("H86").Value = 3220.7144
$ws.Range("I86").Value = 3010
$ws.Range("K86").Value = 3010
$ws.Range("M86").Value = -1887
$ws.Range("H89").Value = 3220.7144
$ws.Range("I89").Value = 3010
$ws.Range("K89").Value = 15050
$ws.Range("M89").Value = -9434
$ws.Range("H94").Value = 5081.8
$ws.Range("I94").Value = 5506
$ws.Range("J94").Value = 4975.75
$ws.Range("K94").Value = 5506
$ws.Range("L94").Value = 4975.75
$ws.Range("M94").Value = -5055
$ws.Range("N94").Value = -5877.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1985.2941
$ws.Range("I117").Value = 929
$ws.Range("J117").Value = 2051.3125
$ws.Range("K117").Value = 2787
$ws.Range("L117").Value = 6153.9375
$ws.Range("M117").Value = 655
$ws.Range("N117").Value = -13037.9375
$ws.Range("H121").Value = 488.16666
$ws.Range("I121").Value = 357.25
$ws.Range("J121").Value = 750
$ws.Range("K121").Value = 1071.75
$ws.Range("L121").Value = 2250
$ws.Range("M121").Value = 238.25
$ws.Range("N121").Value = -4870
$ws.Range("H129").Value = 1454.7407
$ws.Range("I129").Value = 1286.5
$ws.Range("J129").Value = 1502.8096
$ws.Range("K129").Value = 3859.5
$ws.Range("L129").Value = 4508.4288
$ws.Range("M129").Value = 1140.5
$ws.Range("N129").Value = -14508.4288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4843.0713
$ws.Range("I70").Value = 4230.3
$ws.Range("J70").Value = 6375
$ws.Range("K70").Value = 4230.3
$ws.Range("L70").Value = 6375
$ws.Range("M70").Value = -3960.3
$ws.Range("N70").Value = -6915
$ws.Range("H73").Value = 4843.0713
$ws.Range("I73").Value = 4230.3
$ws.Range("J73").Value = 6375
$ws.Range("K73").Value = 4230.3
$ws.Range("L73").Value = 6375
$ws.Range("M73").Value = -3294.3
$ws.Range("N73").Value = -8247
$ws.Range("H80").Value = 3758.75
$ws.Range("I80").Value = 3780.5
$ws.Range("J80").Value = 3650
$ws.Range("K80").Value = 3780.5
$ws.Range("L80").Value = 3650
$ws.Range("M80").Value = -2782.5
$ws.Range("N80").Value = -5646
$ws.Range("H83").Value = 3758.75
$ws.Range("I83").Value = 3780.5
$ws.Range("J83").Value = 3650
$ws.Range("K83").Value = 18902.5
$ws.Range("L83").Value = 18250
$ws.Range("M83").Value = -13910.5
$ws.Range("N83").Value = -28234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 484.17648
$ws.Range("I46").Value = 468.7143
$ws.Range("J46").Value = 495
$ws.Range("K46").Value = 468.7143
$ws.Range("L46").Value = 495
$ws.Range("M46").Value = -280.7143
$ws.Range("N46").Value = -871

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 28127
$ws.Range("J93").Value = 28127
$ws.Range("L93").Value = 28127
$ws.Range("N93").Value = -33119
